$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the TOC sheet with two new rows describing the new tabs.
# ---------------------------------------------------------------------------
$toc = $wb.Worksheets.Item("TOC")
$toc.Range("A19").Value = "Disablity"
$toc.Range("B19").Value = "People with Disabilities (%) by County and SCAG Region"
$toc.Range("A20").Value = "Categories_Disablity"
$toc.Range("B20").Value = "Disabilites by 6 categories (%) and then by County and SCAG Region: Self-care, Hearing, Vision, Independent, Ambulatory, Cognitive"

# ---------------------------------------------------------------------------
# 2) Add the "Disablity" sheet (overall disability % by county).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dis = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$dis.Name = "Disablity"

$dis.Range("A1").Value = "county"
$dis.Range("B1").Value = "dis_perc"
$dis.Range("A1:B1").Font.Bold = $true
$dis.Range("A1:B1").HorizontalAlignment = -4108

$disData = @(
    @("Imperial", 14.44),
    @("Los Angeles", 10.42),
    @("Orange", 8.880000000000001),
    @("Riverside", 11.65),
    @("San Bernardino", 11.77),
    @("Ventura", 11.18),
    @("SCAG", 10.54)
)

$r = 2
foreach ($row in $disData) {
    $dis.Cells.Item($r, 1).Value = $row[0]
    $dis.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Add the "Categories_Disablity" sheet (disability % by category/county).
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$cat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$cat.Name = "Categories_Disablity"

$cat.Range("A1").Value = "county"
$cat.Range("B1").Value = "disability_type"
$cat.Range("C1").Value = "percentage"
$cat.Range("A1:C1").Font.Bold = $true
$cat.Range("A1:C1").HorizontalAlignment = -4108

$catData = @(
    @("Imperial", "Self-care difficulty", 3.89),
    @("Los Angeles", "Self-care difficulty", 3.02),
    @("Orange", "Self-care difficulty", 2.2),
    @("Riverside", "Self-care difficulty", 2.66),
    @("San Bernardino", "Self-care difficulty", 2.67),
    @("Ventura", "Self-care difficulty", 2.56),
    @("SCAG", "Self-care difficulty", 2.78),

    @("Imperial", "Hearing difficulty", 3.01),
    @("Los Angeles", "Hearing difficulty", 2.56),
    @("Orange", "Hearing difficulty", 2.54),
    @("Riverside", "Hearing difficulty", 3.35),
    @("San Bernardino", "Hearing difficulty", 3.19),
    @("Ventura", "Hearing difficulty", 3.32),
    @("SCAG", "Hearing difficulty", 2.77),

    @("Imperial", "Vision difficulty", 2.12),
    @("Los Angeles", "Vision difficulty", 2.03),
    @("Orange", "Vision difficulty", 1.56),
    @("Riverside", "Vision difficulty", 2.25),
    @("San Bernardino", "Vision difficulty", 2.42),
    @("Ventura", "Vision difficulty", 2.03),
    @("SCAG", "Vision difficulty", 2.02),

    @("Imperial", "Independent living", 5.78),
    @("Los Angeles", "Independent living", 4.7),
    @("Orange", "Independent living", 3.65),
    @("Riverside", "Independent living", 4.5),
    @("San Bernardino", "Independent living", 4.55),
    @("Ventura", "Independent living", 4.53),
    @("SCAG", "Independent living", 4.48),

    @("Imperial", "Ambulatory difficulty", 7.94),
    @("Los Angeles", "Ambulatory difficulty", 5.65),
    @("Orange", "Ambulatory difficulty", 4.39),
    @("Riverside", "Ambulatory difficulty", 6.07),
    @("San Bernardino", "Ambulatory difficulty", 5.87),
    @("Ventura", "Ambulatory difficulty", 5.51),
    @("SCAG", "Ambulatory difficulty", 5.53),

    @("Imperial", "Cognitive difficulty", 6.56),
    @("Los Angeles", "Cognitive difficulty", 4.17),
    @("Orange", "Cognitive difficulty", 3.35),
    @("Riverside", "Cognitive difficulty", 4.31),
    @("San Bernardino", "Cognitive difficulty", 4.9),
    @("Ventura", "Cognitive difficulty", 4.27),
    @("SCAG", "Cognitive difficulty", 4.16)
)

$r = 2
foreach ($row in $catData) {
    $cat.Cells.Item($r, 1).Value = $row[0]
    $cat.Cells.Item($r, 2).Value = $row[1]
    $cat.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Re-select the TOC sheet so the active tab matches the original workbook.
$toc.Select()
